$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '30.313.46'
Set-TextValue $ws.Range("E2") '  -0.31%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.084.25'
Set-TextValue $ws.Range("E3") '  +3.29%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.9990'
Set-TextValue $ws.Range("E4") '  -0.32%  '

# Row 5
Set-TextValue $ws.Range("D5") '328.60'
Set-TextValue $ws.Range("E5") '  +1.09%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.9994'
Set-TextValue $ws.Range("E6") '  -0.18%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5224'
Set-TextValue $ws.Range("E7") '  +1.67%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.4320'
Set-TextValue $ws.Range("E8") '  +2.55%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.08837'
Set-TextValue $ws.Range("E9") '  +1.25%  '

# Row 10
Set-TextValue $ws.Range("E10") '  +7.37%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.163'
Set-TextValue $ws.Range("E11") '  +2.36%  '

# Row 12
Set-TextValue $ws.Range("D12") '24.51'
Set-TextValue $ws.Range("E12") '  -1.44%  '

# Row 13
Set-TextValue $ws.Range("D13") '2.082.51'
Set-TextValue $ws.Range("E13") '  +3.43%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.745'
Set-TextValue $ws.Range("E14") '  +2.12%  '

# Row 15
Set-TextValue $ws.Range("D15") '7.699'
Set-TextValue $ws.Range("E15") '  +2.89%  '

# Row 16
Set-TextValue $ws.Range("D16") '95.60'
Set-TextValue $ws.Range("E16") '  +1.05%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.9993'
Set-TextValue $ws.Range("E17") '  -0.30%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.00001126'
Set-TextValue $ws.Range("E18") '  +0.99%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06631'
Set-TextValue $ws.Range("E19") '  +1.56%  '

# Row 20
Set-TextValue $ws.Range("D20") '18.92'
Set-TextValue $ws.Range("E20") '  -0.36%  '

# Row 21
Set-TextValue $ws.Range("D21") '0.9988'
Set-TextValue $ws.Range("E21") '  -0.23%  '

# Row 22
Set-TextValue $ws.Range("D22") '6.323'
Set-TextValue $ws.Range("E22") '  +1.82%  '

# Row 23
Set-TextValue $ws.Range("D23") '30.372.40'
Set-TextValue $ws.Range("E23") '  -0.31%  '

# Row 24
Set-TextValue $ws.Range("D24") '12.37'
Set-TextValue $ws.Range("E24") '  +4.29%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.297'
Set-TextValue $ws.Range("E25") '  +2.52%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.328.85'
Set-TextValue $ws.Range("E26") '  +3.50%  '

# Row 27
Set-TextValue $ws.Range("D27") '22.38'
Set-TextValue $ws.Range("E27") '  -0.20%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.592'
Set-TextValue $ws.Range("E28") '  +6.60%  '

# Row 29
Set-TextValue $ws.Range("D29") '161.88'
Set-TextValue $ws.Range("E29") '  -0.55%  '

# Row 30
Set-TextValue $ws.Range("D30") '131.65'
Set-TextValue $ws.Range("E30") '  +0.04%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.198'
Set-TextValue $ws.Range("E31") '  +5.00%  '

# Row 32
Set-TextValue $ws.Range("E32") '  +1.56%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.666'
Set-TextValue $ws.Range("E33") '  +21.29%  '

# Row 34
Set-TextValue $ws.Range("D34") '6.186'
Set-TextValue $ws.Range("E34") '  +1.83%  '

# Row 35
Set-TextValue $ws.Range("D35") '3.871'
Set-TextValue $ws.Range("E35") '  +1.08%  '

# Row 36
Set-TextValue $ws.Range("D36") '10.01'
Set-TextValue $ws.Range("E36") '  +10.37%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.02572'
Set-TextValue $ws.Range("E37") '  +1.66%  '

# Row 38
Set-TextValue $ws.Range("B38") 'Hedera'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D38") '0.06676'
Set-TextValue $ws.Range("E38") '  +0.12%  '

# Row 39
Set-TextValue $ws.Range("B39") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D39") '5.459'
Set-TextValue $ws.Range("E39") '  -0.44%  '

# Row 40
Set-TextValue $ws.Range("B40") 'Aptos'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D40") '12.70'
Set-TextValue $ws.Range("E40") '  +3.00%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.2266'
Set-TextValue $ws.Range("E41") '  +2.82%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.6835'
Set-TextValue $ws.Range("E42") '  +2.44%  '

# Row 43
Set-TextValue $ws.Range("D43") '1.245'
Set-TextValue $ws.Range("E43") '  +0.76%  '

# Row 44
Set-TextValue $ws.Range("B44") 'Frax'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D44") '0.9983'
Set-TextValue $ws.Range("E44") '  -0.20%  '

# Row 45
Set-TextValue $ws.Range("B45") 'EnergySwap'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '14.07'
Set-TextValue $ws.Range("E45") '  +2.45%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.6390'
Set-TextValue $ws.Range("E46") '  +3.32%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.209'
Set-TextValue $ws.Range("E47") '  +0.73%  '

# Row 48
Set-TextValue $ws.Range("D48") '3.607'
Set-TextValue $ws.Range("E48") '  -1.47%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.250'
Set-TextValue $ws.Range("E49") '  -1.26%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.192'
Set-TextValue $ws.Range("E50") '  +7.39%  '

# Row 51
Set-TextValue $ws.Range("D51") '81.71'
Set-TextValue $ws.Range("E51") '  +0.81%  '
